# Auto-generated Excel COM-interop edit script
# Applies the numeric corrections described in the commit diff,
# sheet by sheet, cell by cell (set new value, or clear cell when removed).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 6533
$ws.Range("J43").Value = 6533
$ws.Range("L43").Value = 6533
$ws.Range("N43").Value = -6671
$ws.Range("H112").Value = 3158.9092
$ws.Range("J112").Value = 3335
$ws.Range("L112").Value = 10005
$ws.Range("N112").Value = -12221
$ws.Range("H137").Value = 1677.8948
$ws.Range("I137").Value = 1440.1666
$ws.Range("K137").Value = 4320.4998
$ws.Range("M137").Value = -1770.4998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("M61").ClearContents()
$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("M136").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 4441.364
$ws.Range("I20").Value = 2650.5557
$ws.Range("K20").Value = 2650.5557
$ws.Range("M20").Value = -2403.5557

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("M16").ClearContents()
$ws.Range("N16").ClearContents()
$ws.Range("H22").Value = 1056.7142
$ws.Range("I22").Value = 1157.8334
$ws.Range("J22").Value = 450
$ws.Range("K22").Value = 1157.8334
$ws.Range("L22").Value = 450
$ws.Range("M22").Value = -807.8334
$ws.Range("N22").Value = -1150
$ws.Range("H32").Value = 3107
$ws.Range("I32").Value = 3107
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 3107
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -2791
$ws.Range("N32").ClearContents()
$ws.Range("H41").Value = 15879.25
$ws.Range("J41").Value = 17425
$ws.Range("L41").Value = 17425
$ws.Range("N41").Value = -18281
$ws.Range("H58").Value = 2527.5
$ws.Range("I58").Value = 1703.3334
$ws.Range("K58").Value = 1703.3334
$ws.Range("M58").Value = -1500.3334
$ws.Range("H74").Value = 33998.5
$ws.Range("J74").Value = 33998.5
$ws.Range("L74").Value = 33998.5
$ws.Range("N74").Value = -35746.5
$ws.Range("H77").Value = 33998.5
$ws.Range("J77").Value = 33998.5
$ws.Range("L77").Value = 101995.5
$ws.Range("N77").Value = -110731.5
$ws.Range("H99").Value = 8379.4
$ws.Range("I99").Value = 7299.6665
$ws.Range("J99").Value = 9999
$ws.Range("K99").Value = 7299.6665
$ws.Range("L99").Value = 9999
$ws.Range("M99").Value = -5801.6665
$ws.Range("N99").Value = -12995
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("N113").ClearContents()
$ws.Range("H122").Value = 600
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").ClearContents()
$ws.Range("H126").Value = 8379.4
$ws.Range("I126").Value = 7299.6665
$ws.Range("J126").Value = 9999
$ws.Range("K126").Value = 21898.9995
$ws.Range("L126").Value = 29997
$ws.Range("M126").Value = -19428.9995
$ws.Range("N126").Value = -34937
$ws.Range("H132").Value = 5307.727
$ws.Range("I132").Value = 5831.6665
$ws.Range("K132").Value = 17494.9995
$ws.Range("M132").Value = -14964.9995
$ws.Range("H136").Value = 2527.5
$ws.Range("I136").Value = 1703.3334
$ws.Range("K136").Value = 5110.0002
$ws.Range("M136").Value = -2560.0002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 653.44446
$ws.Range("J12").Value = 839.2857
$ws.Range("L12").Value = 2517.8571
$ws.Range("N12").Value = -2863.8571

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3075
$ws.Range("J80").Value = 3150
$ws.Range("L80").Value = 3150
$ws.Range("N80").Value = -5146
$ws.Range("H83").Value = 3075
$ws.Range("J83").Value = 3150
$ws.Range("L83").Value = 15750
$ws.Range("N83").Value = -25734
$ws.Range("H97").Value = 1204.375
$ws.Range("I97").Value = 756.6667
$ws.Range("J97").Value = 2547.5
$ws.Range("K97").Value = 756.6667
$ws.Range("L97").Value = 2547.5
$ws.Range("M97").Value = -260.6667
$ws.Range("N97").Value = -3539.5
$ws.Range("H113").Value = 1608.25
$ws.Range("I113").Value = 1744
$ws.Range("K113").Value = 1744
$ws.Range("M113").Value = 426

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1435.15
$ws.Range("I46").Value = 1380
$ws.Range("J46").Value = 1490.3
$ws.Range("K46").Value = 1380
$ws.Range("L46").Value = 1490.3
$ws.Range("M46").Value = -1192
$ws.Range("N46").Value = -1866.3
$ws.Range("H55").Value = 421.33334
$ws.Range("I55").Value = 68.666664
$ws.Range("K55").Value = 68.666664
$ws.Range("M55").Value = 104.333336
$ws.Range("H68").Value = 1900
$ws.Range("I68").Value = 1900
$ws.Range("K68").Value = 1900
$ws.Range("M68").Value = -1151
$ws.Range("H71").Value = 1900
$ws.Range("I71").Value = 1900
$ws.Range("K71").Value = 9500
$ws.Range("M71").Value = -5756
$ws.Range("H93").Value = 2248.5
$ws.Range("I93").Value = 1997.5
$ws.Range("K93").Value = 1997.5
$ws.Range("M93").Value = -749.5
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4166.6665
$ws.Range("I62").Value = 2750
$ws.Range("J62").Value = 7000
$ws.Range("K62").Value = 2750
$ws.Range("L62").Value = 7000
$ws.Range("M62").Value = -2126
$ws.Range("N62").Value = -8248
$ws.Range("H65").Value = 4166.6665
$ws.Range("I65").Value = 2750
$ws.Range("J65").Value = 7000
$ws.Range("K65").Value = 13750
$ws.Range("L65").Value = 35000
$ws.Range("M65").Value = -10630
$ws.Range("N65").Value = -41240
$ws.Range("H136").Value = 2599
$ws.Range("I136").Value = 1549.1666
$ws.Range("K136").Value = 4647.4998
$ws.Range("M136").Value = -2097.4998
